# Update ObjTables document/table attribute names to lowerCamelCase
# (Type -> type, Id -> id, Name -> name, Date -> date, ObjTablesVersion -> objTablesVersion)
# and refresh the date/time stamps, per commit:
#   "changing document, table attributes to lowerCamelCase"

$wb = $excel.ActiveWorkbook

# Document-level banner on the first sheet (row A2) uses the !!! prefix.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Unprotect()
$ws1.Range("A2").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 13:01:14'"
$ws1.Protect($null, $true, $true, $true)

# Table-level headers (row A1) on every data sheet, in workbook order.
$tableUpdates = @(
    @{ Index = 1;  Id = "Compartment";            Date = "2020-03-09 13:01:14" },
    @{ Index = 2;  Id = "Compound";                Date = "2020-03-09 13:01:14" },
    @{ Index = 3;  Id = "Definition";              Date = "2020-03-09 13:01:14" },
    @{ Index = 4;  Id = "Enzyme";                  Date = "2020-03-09 13:01:14" },
    @{ Index = 5;  Id = "FbcObjective";            Date = "2020-03-09 13:01:14" },
    @{ Index = 6;  Id = "Gene";                    Date = "2020-03-09 13:01:14" },
    @{ Index = 7;  Id = "Layout";                  Date = "2020-03-09 13:01:14" },
    @{ Index = 8;  Id = "Measurement";             Date = "2020-03-09 13:01:14" },
    @{ Index = 9;  Id = "PbConfig";                Date = "2020-03-09 13:01:14" },
    @{ Index = 10; Id = "Position";                Date = "2020-03-09 13:01:14" },
    @{ Index = 11; Id = "Protein";                 Date = "2020-03-09 13:01:14" },
    @{ Index = 12; Id = "Quantity";                Date = "2020-03-09 13:01:14" },
    @{ Index = 13; Id = "QuantityInfo";            Date = "2020-03-09 13:01:15" },
    @{ Index = 14; Id = "QuantityMatrix";          Date = "2020-03-09 13:01:15" },
    @{ Index = 15; Id = "Reaction";                Date = "2020-03-09 13:01:15" },
    @{ Index = 16; Id = "ReactionStoichiometry";   Date = "2020-03-09 13:01:15" },
    @{ Index = 17; Id = "Regulator";               Date = "2020-03-09 13:01:15" },
    @{ Index = 18; Id = "Relation";                Date = "2020-03-09 13:01:15" },
    @{ Index = 19; Id = "Relationship";            Date = "2020-03-09 13:01:15" },
    @{ Index = 20; Id = "SparseMatrix";            Date = "2020-03-09 13:01:15" },
    @{ Index = 21; Id = "SparseMatrixColumn";      Date = "2020-03-09 13:01:15" },
    @{ Index = 22; Id = "SparseMatrixOrdered";     Date = "2020-03-09 13:01:15" },
    @{ Index = 23; Id = "SparseMatrixRow";         Date = "2020-03-09 13:01:15" },
    @{ Index = 24; Id = "StoichiometricMatrix";    Date = "2020-03-09 13:01:15" },
    @{ Index = 25; Id = "rxnconContingencyList";   Date = "2020-03-09 13:01:15" },
    @{ Index = 26; Id = "rxnconReactionList";      Date = "2020-03-09 13:01:15" }
)

foreach ($entry in $tableUpdates) {
    $ws = $wb.Worksheets.Item($entry.Index)
    $id = $entry.Id
    $value = "!!ObjTables type='Data' id='$id' name='$id' date='$($entry.Date)' objTablesVersion='0.0.8'"
    $ws.Unprotect()
    $ws.Range("A1").Value = $value
    $ws.Protect($null, $true, $true, $true)
}
